$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "-1-Change The Name of Providers to ProviderProfile."
#       -> "-1-Make it to be matcheable logically between Web and Test Project"
#    (paragraph 3 — no structural changes yet, so indices are still the
#    original ones)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("-1-Change The Name of Providers to ", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "-1-Make it to be ", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("ProviderProfile", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "matcheable", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("matcheable.", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "matcheable logically between Web and Test Project", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Make it to be matchebble logically"  (paragraph 4)
#       -> "-2-User profile Dependency Injection ( Will be adjusted based on the control behavior"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Make it to be matchebble logically", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "-2-User profile Dependency Injection ( Will be adjusted based on the control behavior", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert a brand-new paragraph right after paragraph 4:
#       "Too much to do with UserHelperHtmlExtension and the views)"
#    formatted with a first-line indent of 720 twips (36 pt).
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(5)
$newPara.Range.Text = "Too much to do with UserHelperHtmlExtension and the views)"
$newPara.Format.FirstLineIndent = 36

# ---------------------------------------------------------------------------
# 4) "-2-User profile Dependency Injection" (now paragraph 6, previously
#    paragraph 5) -> "-4-Clean Up upload and other WIP Controllers"
#    Scope the Find to that specific paragraph's Range so the similarly
#    worded text in the paragraph we just wrote above is left untouched.
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6).Range
$p6.Find.Execute("-2-User profile Dependency Injection", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "-4-Clean Up upload and other WIP Controllers", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Drop the now-duplicated "-3-Clean Up Solutions" paragraph (7) and the
#    now-duplicated "-4-Clean Up upload and other WIP Controllers" paragraph
#    (still 7 once the first one is gone) that used to follow it.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(7).Range.Delete() | Out-Null
$d.Paragraphs.Item(7).Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 6) Relocate the hidden "_GoBack" bookmark: it used to sit at the very start
#    of the "-3-Unit Test..." paragraph; it now belongs inside
#    "-5-Update Facebook graph icon" (paragraph 7), between "grap" and
#    "h icon".
# ---------------------------------------------------------------------------
$bmRng = $d.Paragraphs.Item(7).Range
$bmRng.Find.Execute("-5-Update Facebook grap") | Out-Null
$bmRng.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

Write-Output "Done."
